# Updates Jogos_do_Dia_Betfair_Back_Lay_2025-10-13 worksheet:
#  - Refreshes Betfair back/lay odds for the existing fixtures (rows 2-15)
#  - Appends four newly listed fixtures (rows 16-19)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update odds on existing rows (2-15) ---
$ws.Range("F2").Value = 2.66
$ws.Range("G2").Value = 3.45
$ws.Range("H2").Value = 2.6
$ws.Range("I2").Value = 3.35
$ws.Range("J2").Value = 2.86
$ws.Range("K2").Value = 3.7
$ws.Range("L2").Value = 1.34
$ws.Range("N2").Value = 1.8
$ws.Range("O2").Value = 1.35
$ws.Range("P2").Value = 1.68
$ws.Range("Q2").Value = 1.89
$ws.Range("R2").Value = 1.09
$ws.Range("S2").Value = 1.89
$ws.Range("V2").Value = 1.42
$ws.Range("W2").Value = 1.4
$ws.Range("G3").Value = 2.14
$ws.Range("W3").Value = 1.89
$ws.Range("F4").Value = 1.57
$ws.Range("G4").Value = 1.64
$ws.Range("I4").Value = 7.2
$ws.Range("Q4").Value = 1.83
$ws.Range("R4").Value = 1.37
$ws.Range("W4").Value = 2.54
$ws.Range("Z4").Value = 55
$ws.Range("AJ4").Value = 980
$ws.Range("AK4").Value = 980
$ws.Range("AN4").Value = 9.4
$ws.Range("F6").Value = 5.4
$ws.Range("G6").Value = 7.6
$ws.Range("I6").Value = 1.84
$ws.Range("J6").Value = 3.2
$ws.Range("N6").Value = 3.2
$ws.Range("P6").Value = 1.76
$ws.Range("Q6").Value = 2.06
$ws.Range("R6").Value = 1.29
$ws.Range("S6").Value = 3.45
$ws.Range("T6").Value = 1.98
$ws.Range("U6").Value = 1.84
$ws.Range("V6").Value = 2
$ws.Range("W6").Value = 1.15
$ws.Range("X6").Value = 15.5
$ws.Range("Y6").Value = 9
$ws.Range("Z6").Value = 12
$ws.Range("AC6").Value = 10.5
$ws.Range("AD6").Value = 12.5
$ws.Range("AJ6").Value = 210
$ws.Range("AL6").Value = 120
$ws.Range("AM6").Value = 180
$ws.Range("G7").Value = 1.9
$ws.Range("H7").Value = 4.7
$ws.Range("I7").Value = 6.4
$ws.Range("J7").Value = 3.75
$ws.Range("R7").Value = 1.4
$ws.Range("S7").Value = 3.1
$ws.Range("T7").Value = 1.78
$ws.Range("W7").Value = 2.12
$ws.Range("X7").Value = 21
$ws.Range("AA7").Value = 140
$ws.Range("AM7").Value = 120
$ws.Range("G8").Value = 1.98
$ws.Range("J8").Value = 3.15
$ws.Range("K8").Value = 5.5
$ws.Range("N8").Value = 1.47
$ws.Range("P8").Value = 1.47
$ws.Range("Q8").Value = 2.28
$ws.Range("S8").Value = 2.28
$ws.Range("W8").Value = 2.02
$ws.Range("F9").Value = 2.22
$ws.Range("H9").Value = 2.6
$ws.Range("I9").Value = 3.2
$ws.Range("N9").Value = 1.76
$ws.Range("P9").Value = 1.76
$ws.Range("Q9").Value = 2.04
$ws.Range("S9").Value = 2.04
$ws.Range("T9").Value = 1.57
$ws.Range("U9").Value = 1.78
$ws.Range("V9").Value = 1.45
$ws.Range("M10").Value = 1.15
$ws.Range("N10").Value = 2.18
$ws.Range("O10").Value = 1.71
$ws.Range("P10").Value = 1.39
$ws.Range("U10").Value = 1.56
$ws.Range("Y10").Value = 9.800000000000001
$ws.Range("AA10").Value = 140
$ws.Range("AB10").Value = 6
$ws.Range("AD10").Value = 980
$ws.Range("AE10").Value = 110
$ws.Range("AH10").Value = 36
$ws.Range("AI10").Value = 140
$ws.Range("AJ10").Value = 980
$ws.Range("AK10").Value = 40
$ws.Range("AL10").Value = 130
$ws.Range("AM10").Value = 400
$ws.Range("AN10").Value = 44
$ws.Range("L11").Value = 1.55
$ws.Range("N11").Value = 2.6
$ws.Range("AJ11").Value = 1000
$ws.Range("H12").Value = 2.78
$ws.Range("I12").Value = 3.25
$ws.Range("N12").Value = 1.59
$ws.Range("V12").Value = 1.44
$ws.Range("J13").Value = 2.66
$ws.Range("F14").Value = 2.08
$ws.Range("G14").Value = 2.38
$ws.Range("H14").Value = 3.15
$ws.Range("I14").Value = 4.7
$ws.Range("J14").Value = 2.9
$ws.Range("K14").Value = 4.3
$ws.Range("M14").Value = 1.07
$ws.Range("N14").Value = 2.8
$ws.Range("O14").Value = 1.37
$ws.Range("R14").Value = 1.28
$ws.Range("S14").Value = 3.1
$ws.Range("T14").Value = 1.83
$ws.Range("U14").Value = 1.95
$ws.Range("V14").Value = 1.27
$ws.Range("W14").Value = 1.72
$ws.Range("H15").Value = 2.58
$ws.Range("I15").Value = 4.7
$ws.Range("J15").Value = 2.7
$ws.Range("K15").Value = 5
$ws.Range("R15").Value = 1.1
$ws.Range("V15").Value = 1.27
$ws.Range("AO15").Value = 1000

# --- Append new fixture rows (16-19) ---
# Row 16
$ws.Range("A16").Value = "FIFA World Cup Qualifiers - Americas"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "2025-10-13"
$ws.Range("C16").Value = "21:00:00"
$ws.Range("D16").Value = "Honduras"
$ws.Range("E16").Value = "Haiti"
$ws.Range("F16").Value = 1.82
$ws.Range("G16").Value = 2.04
$ws.Range("H16").Value = 3.9
$ws.Range("I16").Value = 5.2
$ws.Range("J16").Value = 3.35
$ws.Range("K16").Value = 4.4
$ws.Range("L16").Value = 1.45
$ws.Range("M16").Value = 1.04
$ws.Range("N16").Value = 1.04
$ws.Range("O16").Value = 1.36
$ws.Range("P16").Value = 1.24
$ws.Range("Q16").Value = 1.36
$ws.Range("R16").Value = 1.18
$ws.Range("S16").Value = 1.37
$ws.Range("T16").Value = 1.01
$ws.Range("U16").Value = 1.01
$ws.Range("V16").Value = 1.25
$ws.Range("W16").Value = 1.96
$ws.Range("X16").Value = 14.5
$ws.Range("Y16").Value = 19.5
$ws.Range("Z16").Value = 42
$ws.Range("AA16").Value = 150
$ws.Range("AB16").Value = 10
$ws.Range("AC16").Value = 9.6
$ws.Range("AD16").Value = 23
$ws.Range("AE16").Value = 85
$ws.Range("AF16").Value = 12.5
$ws.Range("AG16").Value = 12
$ws.Range("AH16").Value = 26
$ws.Range("AI16").Value = 100
$ws.Range("AJ16").Value = 23
$ws.Range("AK16").Value = 24
$ws.Range("AL16").Value = 50
$ws.Range("AM16").Value = 160
$ws.Range("AN16").Value = 18
$ws.Range("AO16").Value = 110

# Row 17
$ws.Range("A17").Value = "Brazilian Serie B"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "2025-10-13"
$ws.Range("C17").Value = "21:30:00"
$ws.Range("D17").Value = "CRB"
$ws.Range("E17").Value = "Ferroviaria"
$ws.Range("F17").Value = 1.62
$ws.Range("G17").Value = 1.73
$ws.Range("H17").Value = 6.2
$ws.Range("I17").Value = 8.199999999999999
$ws.Range("J17").Value = 3.65
$ws.Range("K17").Value = 4.1
$ws.Range("L17").Value = 1.43
$ws.Range("M17").Value = 1.07
$ws.Range("N17").Value = 3.2
$ws.Range("O17").Value = 1.38
$ws.Range("P17").Value = 1.76
$ws.Range("Q17").Value = 2.1
$ws.Range("R17").Value = 1.29
$ws.Range("S17").Value = 3.85
$ws.Range("T17").Value = 2.06
$ws.Range("U17").Value = 1.82
$ws.Range("V17").Value = 1.16
$ws.Range("W17").Value = 2.36
$ws.Range("X17").Value = 15
$ws.Range("Y17").Value = 1000
$ws.Range("Z17").Value = 1000
$ws.Range("AA17").Value = 250
$ws.Range("AB17").Value = 7.2
$ws.Range("AC17").Value = 9
$ws.Range("AD17").Value = 30
$ws.Range("AE17").Value = 140
$ws.Range("AF17").Value = 9.6
$ws.Range("AG17").Value = 10.5
$ws.Range("AH17").Value = 30
$ws.Range("AI17").Value = 130
$ws.Range("AJ17").Value = 17
$ws.Range("AK17").Value = 1000
$ws.Range("AL17").Value = 55
$ws.Range("AM17").Value = 200
$ws.Range("AN17").Value = 12.5
$ws.Range("AO17").Value = 1000

# Row 18
$ws.Range("A18").Value = "Colombian Primera A"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "2025-10-13"
$ws.Range("C18").Value = "21:30:00"
$ws.Range("D18").Value = "Fortaleza FC"
$ws.Range("E18").Value = "Boyaca Chico"
$ws.Range("F18").Value = 1.83
$ws.Range("G18").Value = 2.06
$ws.Range("H18").Value = 4.8
$ws.Range("I18").Value = 6.4
$ws.Range("J18").Value = 3.15
$ws.Range("K18").Value = 3.7
$ws.Range("L18").Value = 1.01
$ws.Range("M18").Value = 1.09
$ws.Range("N18").Value = 2.72
$ws.Range("O18").Value = 1.46
$ws.Range("P18").Value = 1.58
$ws.Range("Q18").Value = 2.18
$ws.Range("R18").Value = 1.21
$ws.Range("S18").Value = 4.2
$ws.Range("T18").Value = 2.08
$ws.Range("U18").Value = 1.74
$ws.Range("V18").Value = 1.2
$ws.Range("W18").Value = 1.99
$ws.Range("X18").Value = 12
$ws.Range("Y18").Value = 15.5
$ws.Range("Z18").Value = 1000
$ws.Range("AA18").Value = 190
$ws.Range("AB18").Value = 6.8
$ws.Range("AC18").Value = 8.4
$ws.Range("AD18").Value = 23
$ws.Range("AE18").Value = 120
$ws.Range("AF18").Value = 10.5
$ws.Range("AG18").Value = 11.5
$ws.Range("AH18").Value = 27
$ws.Range("AI18").Value = 130
$ws.Range("AJ18").Value = 24
$ws.Range("AK18").Value = 26
$ws.Range("AL18").Value = 1000
$ws.Range("AM18").Value = 230
$ws.Range("AN18").Value = 24
$ws.Range("AO18").Value = 1000

# Row 19
$ws.Range("A19").Value = "FIFA World Cup Qualifiers - Americas"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "2025-10-13"
$ws.Range("C19").Value = "23:00:00"
$ws.Range("D19").Value = "Costa Rica"
$ws.Range("E19").Value = "Nicaragua"
$ws.Range("F19").Value = 1.04
$ws.Range("G19").Value = 1.3
$ws.Range("H19").Value = 4.3
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 4.4
$ws.Range("K19").Value = 950
$ws.Range("L19").Value = 1.34
$ws.Range("M19").Value = 1.01
$ws.Range("N19").Value = 1.03
$ws.Range("O19").Value = 1.24
$ws.Range("P19").Value = 1.89
$ws.Range("Q19").Value = 1.25
$ws.Range("R19").Value = 1.18
$ws.Range("S19").Value = 1.5
$ws.Range("T19").Value = 1.01
$ws.Range("U19").Value = 1.01
$ws.Range("V19").Value = 1.01
$ws.Range("W19").Value = 4.3
$ws.Range("X19").Value = 1000
$ws.Range("Y19").Value = 1000
$ws.Range("Z19").Value = 1000
$ws.Range("AA19").Value = 1000
$ws.Range("AB19").Value = 1000
$ws.Range("AC19").Value = 1000
$ws.Range("AD19").Value = 1000
$ws.Range("AE19").Value = 1000
$ws.Range("AF19").Value = 1000
$ws.Range("AG19").Value = 1000
$ws.Range("AH19").Value = 1000
$ws.Range("AI19").Value = 1000
$ws.Range("AJ19").Value = 1000
$ws.Range("AK19").Value = 1000
$ws.Range("AL19").Value = 1000
$ws.Range("AM19").Value = 1000
$ws.Range("AN19").Value = 1000
$ws.Range("AO19").Value = 1000
